# Re-shuffle the per-record data (observation identity, quantity, coordinates,
# accuracy and start/end time) among rows 59-77 of the "Artfynd" sheet.
#
# The location/date/observer columns (P, T, U, V, W, Y, AA, AW, AX, ...) stay
# put on their row; only the "record" columns below move between rows
# according to the mapping captured from the source diff.
#
# $targetRow -> $sourceRow means: the record values that currently live on
# $sourceRow should end up on $targetRow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that travel together as one "record".
$recordCols = @("A","B","D","E","F","G","H","I","M","Q","R","S","Z","AB")

# Columns holding plain numbers (everything else in $recordCols is text).
$numericCols = @("A","B","E","Q","R","S")

$firstRow = 59
$lastRow = 77

# target row -> source row (where its new content currently lives)
$rowMap = @{
    59 = 66; 60 = 62; 61 = 67; 62 = 61; 63 = 68
    64 = 65; 65 = 76; 66 = 70; 67 = 59; 68 = 63
    69 = 74; 70 = 75; 71 = 64; 72 = 69; 73 = 77
    74 = 72; 75 = 73; 76 = 71; 77 = 60
}

# 1) Snapshot every record's current values before any writes happen, so
#    overlapping/cyclic moves don't clobber values we still need to read.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($col in $recordCols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write each target row's record using the snapshotted source row values.
for ($targetRow = $firstRow; $targetRow -le $lastRow; $targetRow++) {
    $sourceRow = $rowMap[$targetRow]
    $sourceVals = $snapshot[$sourceRow]
    foreach ($col in $recordCols) {
        $val = $sourceVals[$col]
        $cell = $ws.Range("$col$targetRow")
        if ($null -eq $val -or $val -eq "") {
            $cell.Value2 = ""
        } elseif ($numericCols -contains $col) {
            $cell.Value2 = [double]$val
        } else {
            $cell.Value2 = [string]$val
        }
    }
}
